# Update the teacher's weekly schedule (Euclides) with the new class preferences.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (7:50)
$ws.Range("B3").Value = "Gestão - MCT-1A"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Gestão - MEC-1A"
$ws.Range("F3").Value = "-"

# Row 4 (8:40)
$ws.Range("C4").Value = "Mecanica material - MEC-2A"

# Row 6 (9:50)
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "Mecanica material - MEC-2A"
$ws.Range("D6").Value = "Usinagem - MCT-3A"
$ws.Range("F6").Value = "Usinagem - MEC-3A"
